$wb = $excel.ActiveWorkbook

# --- Add the new "gen_res|pmax" sheet, placed after "gen|pmax" ---
$genPmax = $wb.Worksheets.Item("gen|pmax")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $genPmax)
$newSheet.Name = "gen_res|pmax"

$newSheet.Range("A1").Value = "Time\Id"
$newSheet.Range("B1").Value = 1
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 1000
$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = 1000
$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = 500
[void]$newSheet.Range("B5").Select()

# --- Update the "gen|pmax" sheet data: B4 6000 -> 5000 ---
$wb.Worksheets.Item("gen|pmax").Range("B4").Value = 5000

# --- Insert a row for "gen_res|pmax" in the ReadMe sheet, right before "storage|inflow" ---
$readMe = $wb.Worksheets.Item("ReadMe")
$readMe.Rows.Item(10).Insert()
$readMe.Rows.Item(10).RowHeight = 30

$readMe.Range("A10").Value = "gen_res|pmax"
$readMe.Range("B10").Formula = "=COUNT('gen_res|pmax'!`$1:`$1)"
$readMe.Range("C10").Formula = "=IF(B10,COUNT('gen_res|pmax'!B:B)-1,0)"
$readMe.Range("D10").Value = "MW"
$readMe.Range("E10").Value = "Available production. If MW: in MW.`nIf %: 0.5 means 50% of the production rating."

[void]$readMe.Range("A11").Select()
[void]$readMe.Activate()
